$d = $word.ActiveDocument
$p31 = $d.Paragraphs.Item(31)
$rng = $p31.Range.Duplicate
$rng.Collapse(1) # wdCollapseStart
Write-Host "start=$($rng.Start) end=$($rng.End)"
$d.Bookmarks.Add("TEST_P31_START", $rng)
